# Deploying to gh-pages from @ LinuxForHealth/alvearie-fhir-ig
# Rebrand from "Alvearie"/ibm.com to "LinuxForHealth" + bump to the
# regenerated FHIR Implementation Guide's StructureDefinition spreadsheet
# (CarveOut extension), 7.0.0 -> 8.0.0.

$wb = $excel.ActiveWorkbook

# --- "Metadata" sheet: top-level StructureDefinition properties ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/carve-out"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- "Elements" sheet: per-element definition table ---
$elements = $wb.Worksheets.Item("Elements")
# Extension (row 2): Constraint(s) no longer echoes the inherited ele-1/ext-1
# text on the root slice -- regenerated IG leaves it blank here.
$elements.Range("AI2").Value = ""
# Extension.url (row 5): Fixed Value follows the new canonical URL.
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/carve-out"
